# Add a new column O ("Column with NULL and then mixed") to the TestData
# sheet: a header in row 3, a blank (NULL) cell in row 4, a float in row 5,
# a string in row 6, and a blank cell in row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = "Column with NULL and then mixed"
$ws.Range("O5").Value = 3.4
$ws.Range("O6").Value = "HKEJW"

# Size the new column the way Excel would after typing the header text in
# (best-fit width for the new content).
$ws.Columns.Item(15).ColumnWidth = 29.58203125

# Leave the same kind of post-edit selection state captured in the workbook.
$ws.Range("O7").Select()
